$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 11 (row 11 "CannotNavigateToDoPage" test moved conceptually;
#    actually: pushes the existing empty spacer row 11 down to row 12, and everything below shifts
#    down by one, matching the diff's row renumbering).
$ws.Rows("11:11").Insert()

# 2. Row 9 ("Logout" test case) - rename test + reorder description/outcome strings
$ws.Range("B9").Value2 = "test_<Logout>"
$ws.Range("C9").Value2 = "This is to test that users can logout"
$ws.Range("E9").Value2 = "Message saying ""Logged out"" is shown"

# 3. Row 10 ("CannotNavigateToDoPage" test case) - renumber + rename test + reorder strings
$ws.Range("A10").Value2 = 7
$ws.Range("B10").Value2 = "test_<CannotNavigateToDoPage>"
$ws.Range("C10").Value2 = "This is to test that users cannot navigate to To Do page without logging in"
$ws.Range("E10").Value2 = "Error message saying ""Please login to see this page."" is shown"

# 4. Row 11 (NEW "CannotNavigateToDoHistory" failing test case)
$ws.Range("A11").Value2 = 8
$ws.Range("B11").Value2 = "test_<CannotNavigateToDoHistory>"
$ws.Range("C11").Value2 = "This is to test that users cannot navigate to To Do History page without logging in"
$ws.Range("D11").Value2 = "NIL"
$ws.Range("E11").Value2 = "Error message saying ""Please login to see this page."" is shown"

# 5. Fix up the formatting on the new row 11 cells that the plain Insert() didn't already give us
#    (it copies the row-above's formats by default, but A/F/G should look like row 9's, and B
#    should look like row 9's "test name" cell too).
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("F9").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 6. Leave the selection on B11, matching where the edit was made
$ws.Range("B11").Select()
